# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape, per commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 4914
    5  = 806
    7  = 1285
    10 = 216
    12 = 8
    14 = 120
    15 = 4333
    16 = 6616
    21 = 51
    22 = 4067
    23 = 433
    24 = 62
    29 = 156
    30 = 330
    31 = 339
    34 = 26
    35 = 1602
    37 = 56
    43 = 84
    44 = 613
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
